$d = $word.ActiveDocument

# Add the new "code-keyword" character style used by the inserted code snippet.
$codeKeywordStyle = $d.Styles.Add("code-keyword", 2)
$codeKeywordStyle.BaseStyle = $d.Styles("DefaultParagraphFont")

# The document currently ends with a single paragraph that only holds the
# _GoBack bookmark. Replace it with the new "First create outcome" text,
# the highlighted Python code sample, and put the bookmark back at the
# end of the last code paragraph, followed by one trailing empty paragraph.
$lastParagraph = $d.Paragraphs($d.Paragraphs.Count)
$targetRange = $lastParagraph.Range
$newContentXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p><w:r><w:t xml:space="preserve">First create outcome </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="code-keyword"/><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>def</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr><w:t>ButtonHide</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr><w:t xml:space="preserve">self, state): </w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="code-keyword"/><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>if</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr><w:t xml:space="preserve"> state == </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr><w:t>QtCore.Qt.Checked</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr><w:t>qbtn.setEnabled</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rStyle w:val="code-keyword"/><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>True</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr><w:t xml:space="preserve">) </w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="code-keyword"/><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>else</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr><w:t>qbtn.setEnabled</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rStyle w:val="code-keyword"/><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>False</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FBEDBB"/></w:rPr><w:t>)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/>'
[void]$targetRange.InsertXML($newContentXml)
